# Update cryptocurrency price (D) and hourly volume/change (E) columns
# to reflect the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.788.83'
$ws.Range("E2").Value = '  +8.16%  '

$ws.Range("D3").Value = '1.775.28'
$ws.Range("E3").Value = '  +4.24%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '''225.25'
$ws.Range("E5").Value = '  +1.75%  '

$ws.Range("D6").Value = '''0.558'
$ws.Range("E6").Value = '  +4.43%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").Value = '''30.74'
$ws.Range("E8").Value = '  +3.08%  '

$ws.Range("D9").Value = '''46.48'
$ws.Range("E9").Value = '  +2.80%  '

$ws.Range("E10").Value = '  +3.71%  '

$ws.Range("D11").Value = '''0.0661'
$ws.Range("E11").Value = '  +3.14%  '

$ws.Range("E12").Value = '  +1.36%  '

$ws.Range("D13").Value = '2.030.86'
$ws.Range("E13").Value = '  +4.38%  '

$ws.Range("D14").Value = '1.775.83'
$ws.Range("E14").Value = '  +4.30%  '

$ws.Range("D15").Value = '''0.628'
$ws.Range("E15").Value = '  +2.07%  '

$ws.Range("D16").Value = '33.782.05'
$ws.Range("E16").Value = '  +8.22%  '

$ws.Range("E17").Value = '  -3.37%  '

$ws.Range("D18").Value = '''4.18'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("D19").Value = '''68.48'
$ws.Range("E19").Value = '  +2.03%  '

$ws.Range("D20").Value = '''251.79'
$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("D21").Value = '0.0₃0738'
$ws.Range("E21").Value = '  +2.25%  '

$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '''10.28'
$ws.Range("E23").Value = '  +1.26%  '

$ws.Range("E24").Value = '  -2.58%  '

$ws.Range("E25").Value = '  -0.91%  '

$ws.Range("D26").Value = '''159.00'
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("D27").Value = '''16.49'
$ws.Range("E27").Value = '  +3.10%  '

$ws.Range("D28").Value = '''0.114'
$ws.Range("E28").Value = '  +1.28%  '

$ws.Range("E29").Value = '  +2.84%  '

$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").Value = '  +4.35%  '

$ws.Range("E32").Value = '  +2.08%  '

$ws.Range("E33").Value = '  +3.26%  '

$ws.Range("E34").Value = '  +5.05%  '

$ws.Range("E35").Value = '  +4.27%  '

$ws.Range("D36").Value = '1.484.51'
$ws.Range("E36").Value = '  -2.59%  '

$ws.Range("D37").Value = '''1.07'
$ws.Range("E37").Value = '  +3.30%  '

$ws.Range("D38").Value = '''0.634'
$ws.Range("E38").Value = '  +2.87%  '

$ws.Range("E39").Value = '  +2.68%  '

$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").Value = '''2.34'
$ws.Range("E41").Value = '  +2.04%  '

$ws.Range("D42").Value = '''2.69'
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("E43").Value = '  +3.89%  '

$ws.Range("E44").Value = '  +2.31%  '

$ws.Range("E45").Value = '  +1.92%  '

$ws.Range("E46").Value = '  +3.32%  '

$ws.Range("D47").Value = '1.929.94'
$ws.Range("E47").Value = '  +5.44%  '

$ws.Range("E48").Value = '  +3.43%  '

$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("E50").Value = '  +13.59%  '

$ws.Range("D51").Value = '''50.73'
$ws.Range("E51").Value = '  -2.99%  '
